$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.539.80"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.917.99"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.40"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4822"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2896"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06702"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.34"
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.92"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").Value = "1.919.42"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07558"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.279"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6674"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "294.33"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "30.543.42"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007596"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.546"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "2.169.77"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.433"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.450"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.95"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1064"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("E30").Value = "  +5.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.132"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04998"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7394"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02016"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.680"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.011"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4413"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8659"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.824"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.205"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.45"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.226"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1230"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2528"
$ws.Range("E51").Value = "  +0.15%  "
